$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename shared string "old" -> "Top-Left Symbol" (cell A1 keeps using it) ---
$ws.Range("A1").Value = "Top-Left Symbol"

# --- 2. New column K: header + "copy&paste" label ---
$ws.Range("K2").Value = "copy&paste"

# --- 3. New rows 28-41: same A-D constants as the rest of the table, E=2, F increasing ---
$ws.Range("A28:A41").Value = 1.069
$ws.Range("B28:B41").Value = 0.533
$ws.Range("C28:C41").Value = 2.026
$ws.Range("D28:D41").Value = 1.494
$ws.Range("E28:E41").Value = 2

$ws.Cells.Item(28, 6).Value = 1
for ($r = 29; $r -le 41; $r++) {
    $ws.Cells.Item($r, 6).Value = $r - 27
}

# --- 4. G/H/I/J formulas for the new rows ---
# Row 28 continues right after the existing 25:27 block (kept as-is / not re-touched).
$ws.Range("G28").Formula = "=A28+(E28*1.173)"
$ws.Range("H28").Formula = "=B28+(F28*1.28)"
$ws.Range("I28").Formula = "=C28+(E28*1.173)"
$ws.Range("J28").Formula = "=D28+(F28*1.28)"

# Rows 29:41 form a fresh shared-formula block.
$ws.Range("G29:G41").Formula = "=A29+(E29*1.173)"
$ws.Range("H29:H41").Formula = "=B29+(F29*1.28)"
$ws.Range("I29:I41").Formula = "=C29+(E29*1.173)"
$ws.Range("J29:J41").Formula = "=D29+(F29*1.28)"

# --- 5. New column K: "copy&paste" concatenation of G,H,I,J for every data row ---
$ws.Range("K5").Formula = "=CONCATENATE(G5,"", "",H5,"", "",I5,"", "",J5)"
$ws.Range("K6:K41").Formula = "=CONCATENATE(G6,"", "",H6,"", "",I6,"", "",J6)"

# --- 6. Column K width ---
$ws.Columns.Item(11).ColumnWidth = 22.5

# --- 7. Selection matches the post-edit state ---
$ws.Range("K5:K41").Select()
